$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 6).Value = 2.04
$ws.Cells.Item(3, 7).Value = 2.14
$ws.Cells.Item(3, 8).Value = 4.4
$ws.Cells.Item(3, 9).Value = 4.8
$ws.Cells.Item(3, 11).Value = 3.35
$ws.Cells.Item(3, 12).Value = 1.59
$ws.Cells.Item(3, 14).Value = 2.54
$ws.Cells.Item(3, 15).Value = 1.57
$ws.Cells.Item(3, 16).Value = 1.51
$ws.Cells.Item(3, 17).Value = 2.72
$ws.Cells.Item(3, 20).Value = 2.3
$ws.Cells.Item(3, 21).Value = 1.68
$ws.Cells.Item(3, 23).Value = 1.7
$ws.Cells.Item(3, 29).Value = 7.6
$ws.Cells.Item(3, 31).Value = 100
$ws.Cells.Item(3, 34).Value = 980
$ws.Cells.Item(3, 36).Value = 27
$ws.Cells.Item(3, 38).Value = 70
$ws.Cells.Item(3, 39).Value = 250
$ws.Cells.Item(4, 6).Value = 1.72
$ws.Cells.Item(4, 9).Value = 5.5
$ws.Cells.Item(4, 10).Value = 1.27
$ws.Cells.Item(4, 15).Value = 1.25
$ws.Cells.Item(4, 17).Value = 1.25
$ws.Cells.Item(4, 19).Value = 1.26
$ws.Cells.Item(5, 7).Value = 1.86
$ws.Cells.Item(5, 11).Value = 3.8
$ws.Cells.Item(5, 21).Value = 1.7
$ws.Cells.Item(5, 23).Value = 2.16
$ws.Cells.Item(5, 29).Value = 8.4
$ws.Cells.Item(5, 31).Value = 110
$ws.Cells.Item(5, 34).Value = 27
$ws.Cells.Item(6, 6).Value = 3.35
$ws.Cells.Item(6, 7).Value = 4.6
$ws.Cells.Item(6, 8).Value = 2.06
$ws.Cells.Item(6, 9).Value = 2.3
$ws.Cells.Item(6, 11).Value = 4.3
$ws.Cells.Item(6, 14).Value = 3.25
$ws.Cells.Item(6, 16).Value = 1.78
$ws.Cells.Item(6, 21).Value = 1.98
$ws.Cells.Item(6, 22).Value = 1.78
$ws.Cells.Item(6, 23).Value = 1.31
$ws.Cells.Item(6, 36).Value = 95
$ws.Cells.Item(7, 6).Value = 9.4
$ws.Cells.Item(7, 7).Value = 12.5
$ws.Cells.Item(7, 8).Value = 1.31
$ws.Cells.Item(7, 9).Value = 1.37
$ws.Cells.Item(7, 10).Value = 6
$ws.Cells.Item(7, 11).Value = 6.8
$ws.Cells.Item(7, 16).Value = 2.72
$ws.Cells.Item(7, 17).Value = 1.47
$ws.Cells.Item(7, 21).Value = 2
$ws.Cells.Item(7, 23).Value = 1.09
$ws.Cells.Item(7, 24).Value = 30
$ws.Cells.Item(7, 27).Value = 12
$ws.Cells.Item(7, 29).Value = 14.5
$ws.Cells.Item(7, 30).Value = 11
$ws.Cells.Item(7, 31).Value = 13.5
$ws.Cells.Item(7, 33).Value = 42
$ws.Cells.Item(7, 34).Value = 27
$ws.Cells.Item(7, 35).Value = 32
$ws.Cells.Item(7, 36).Value = 440
$ws.Cells.Item(7, 37).Value = 170
$ws.Cells.Item(7, 40).Value = 1000
$ws.Cells.Item(8, 6).Value = 2
$ws.Cells.Item(8, 7).Value = 2.16
$ws.Cells.Item(8, 9).Value = 5.2
$ws.Cells.Item(8, 10).Value = 3.25
$ws.Cells.Item(8, 13).Value = 1.01
$ws.Cells.Item(8, 22).Value = 1.27
$ws.Cells.Item(8, 23).Value = 1.86
$ws.Cells.Item(9, 6).Value = 2.26
$ws.Cells.Item(9, 8).Value = 2.98
$ws.Cells.Item(9, 9).Value = 3.8
$ws.Cells.Item(9, 11).Value = 4.5
$ws.Cells.Item(9, 22).Value = 1.39
$ws.Cells.Item(9, 39).Value = 85
$ws.Cells.Item(10, 8).Value = 2.4
$ws.Cells.Item(10, 9).Value = 2.64
$ws.Cells.Item(10, 10).Value = 3.45
$ws.Cells.Item(10, 11).Value = 3.75
$ws.Cells.Item(10, 12).Value = 1.39
$ws.Cells.Item(10, 14).Value = 3.5
$ws.Cells.Item(10, 15).Value = 1.33
$ws.Cells.Item(10, 16).Value = 1.87
$ws.Cells.Item(10, 17).Value = 1.96
$ws.Cells.Item(10, 20).Value = 1.74
$ws.Cells.Item(10, 24).Value = 15.5
$ws.Cells.Item(10, 25).Value = 12.5
$ws.Cells.Item(10, 26).Value = 20
$ws.Cells.Item(10, 28).Value = 13
$ws.Cells.Item(10, 29).Value = 9.4
$ws.Cells.Item(10, 30).Value = 13.5
$ws.Cells.Item(10, 32).Value = 22
$ws.Cells.Item(10, 34).Value = 20
$ws.Cells.Item(10, 36).Value = 55
$ws.Cells.Item(11, 6).Value = 2.6
$ws.Cells.Item(11, 8).Value = 2.96
$ws.Cells.Item(11, 10).Value = 2.96
$ws.Cells.Item(11, 11).Value = 3.45
$ws.Cells.Item(11, 13).Value = 1.09
$ws.Cells.Item(11, 14).Value = 2.82
$ws.Cells.Item(11, 16).Value = 1.63
$ws.Cells.Item(11, 18).Value = 1.24
$ws.Cells.Item(11, 21).Value = 1.9
$ws.Cells.Item(11, 23).Value = 1.52
$ws.Cells.Item(12, 6).Value = 4.3
$ws.Cells.Item(12, 7).Value = 4.8
$ws.Cells.Item(12, 8).Value = 2.02
$ws.Cells.Item(12, 9).Value = 2.12
$ws.Cells.Item(12, 10).Value = 3.35
$ws.Cells.Item(12, 11).Value = 3.55
$ws.Cells.Item(12, 12).Value = 1.41
$ws.Cells.Item(12, 15).Value = 1.42
$ws.Cells.Item(12, 16).Value = 1.64
$ws.Cells.Item(12, 17).Value = 2.22
$ws.Cells.Item(12, 19).Value = 4.4
$ws.Cells.Item(12, 20).Value = 1.99
$ws.Cells.Item(12, 22).Value = 1.89
$ws.Cells.Item(12, 23).Value = 1.26
$ws.Cells.Item(12, 24).Value = 11.5
$ws.Cells.Item(12, 26).Value = 12
$ws.Cells.Item(12, 27).Value = 32
$ws.Cells.Item(12, 33).Value = 19
$ws.Cells.Item(12, 35).Value = 980
$ws.Cells.Item(12, 37).Value = 75
$ws.Cells.Item(12, 38).Value = 100
$ws.Cells.Item(12, 39).Value = 170
$ws.Cells.Item(13, 6).Value = 1.64
$ws.Cells.Item(13, 7).Value = 1.76
$ws.Cells.Item(13, 8).Value = 6.2
$ws.Cells.Item(13, 9).Value = 8.4
$ws.Cells.Item(13, 11).Value = 4.1
$ws.Cells.Item(13, 13).Value = 1.09
$ws.Cells.Item(13, 14).Value = 2.9
$ws.Cells.Item(13, 15).Value = 1.39
$ws.Cells.Item(13, 17).Value = 2.16
$ws.Cells.Item(13, 18).Value = 1.26
$ws.Cells.Item(13, 20).Value = 2.1
$ws.Cells.Item(13, 21).Value = 1.75
$ws.Cells.Item(13, 22).Value = 1.13
$ws.Cells.Item(13, 23).Value = 2.3
$ws.Cells.Item(13, 32).Value = 1000
$ws.Cells.Item(14, 6).Value = 2.16
$ws.Cells.Item(14, 7).Value = 2.4
$ws.Cells.Item(14, 8).Value = 3.65
$ws.Cells.Item(14, 9).Value = 4.2
$ws.Cells.Item(14, 10).Value = 3.05
$ws.Cells.Item(14, 11).Value = 3.45
$ws.Cells.Item(14, 12).Value = 1.42
$ws.Cells.Item(14, 13).Value = 1.08
$ws.Cells.Item(14, 14).Value = 2.84
$ws.Cells.Item(14, 15).Value = 1.43
$ws.Cells.Item(14, 19).Value = 4.3
$ws.Cells.Item(14, 20).Value = 1.94
$ws.Cells.Item(14, 22).Value = 1.31
$ws.Cells.Item(14, 23).Value = 1.71
$ws.Cells.Item(14, 24).Value = 11
$ws.Cells.Item(14, 25).Value = 12.5
$ws.Cells.Item(14, 27).Value = 110
$ws.Cells.Item(14, 30).Value = 20
$ws.Cells.Item(14, 31).Value = 60
$ws.Cells.Item(14, 32).Value = 14
$ws.Cells.Item(14, 34).Value = 22
$ws.Cells.Item(14, 35).Value = 75
$ws.Cells.Item(14, 38).Value = 980
$ws.Cells.Item(14, 39).Value = 170
$ws.Cells.Item(14, 41).Value = 90
$ws.Cells.Item(15, 6).Value = 1.77
$ws.Cells.Item(15, 7).Value = 1.99
$ws.Cells.Item(15, 8).Value = 4.7
$ws.Cells.Item(15, 9).Value = 5.7
$ws.Cells.Item(15, 10).Value = 3.45
$ws.Cells.Item(15, 11).Value = 4.1
$ws.Cells.Item(15, 12).Value = 1.35
$ws.Cells.Item(15, 14).Value = 3.35
$ws.Cells.Item(15, 15).Value = 1.34
$ws.Cells.Item(15, 17).Value = 2
$ws.Cells.Item(15, 18).Value = 1.31
$ws.Cells.Item(15, 20).Value = 1.89
$ws.Cells.Item(15, 21).Value = 1.92
$ws.Cells.Item(15, 22).Value = 1.23
$ws.Cells.Item(15, 23).Value = 2
$ws.Cells.Item(15, 25).Value = 18.5
$ws.Cells.Item(15, 26).Value = 46
$ws.Cells.Item(15, 27).Value = 150
$ws.Cells.Item(15, 28).Value = 9.6
$ws.Cells.Item(15, 30).Value = 21
$ws.Cells.Item(15, 32).Value = 11.5
$ws.Cells.Item(15, 33).Value = 11
$ws.Cells.Item(15, 34).Value = 21
$ws.Cells.Item(15, 35).Value = 95
$ws.Cells.Item(15, 36).Value = 22
$ws.Cells.Item(15, 37).Value = 22
$ws.Cells.Item(15, 38).Value = 48
$ws.Cells.Item(15, 40).Value = 16.5
$ws.Cells.Item(16, 6).Value = 4.4
$ws.Cells.Item(16, 7).Value = 5.1
$ws.Cells.Item(16, 11).Value = 3.9
$ws.Cells.Item(16, 14).Value = 2.94
$ws.Cells.Item(16, 23).Value = 1.24
$ws.Cells.Item(16, 38).Value = 110
$ws.Cells.Item(17, 6).Value = 3.6
$ws.Cells.Item(17, 8).Value = 2.18
$ws.Cells.Item(17, 9).Value = 2.3
$ws.Cells.Item(17, 10).Value = 3.45
$ws.Cells.Item(17, 11).Value = 3.65
$ws.Cells.Item(17, 17).Value = 2.06
$ws.Cells.Item(17, 22).Value = 1.76
$ws.Cells.Item(17, 24).Value = 14.5
$ws.Cells.Item(17, 29).Value = 7.8
$ws.Cells.Item(17, 30).Value = 11
$ws.Cells.Item(17, 31).Value = 24
$ws.Cells.Item(17, 32).Value = 27
$ws.Cells.Item(17, 33).Value = 16
$ws.Cells.Item(17, 36).Value = 75
$ws.Cells.Item(17, 40).Value = 50
$ws.Cells.Item(17, 41).Value = 19
$ws.Cells.Item(18, 6).Value = 7
$ws.Cells.Item(18, 7).Value = 9.6
$ws.Cells.Item(18, 8).Value = 1.49
$ws.Cells.Item(18, 9).Value = 1.58
$ws.Cells.Item(18, 11).Value = 5.1
$ws.Cells.Item(18, 14).Value = 2.98
$ws.Cells.Item(18, 20).Value = 2.14
$ws.Cells.Item(18, 22).Value = 2.72
$ws.Cells.Item(18, 23).Value = 1.12
$ws.Cells.Item(18, 26).Value = 980
